$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Attributes")

function Set-AttributeRow($row, $id, $name, $valueType, $code, $codeSystem) {
    $ws1.Cells.Item($row, 1).Value = $id
    $ws1.Cells.Item($row, 3).Value = $valueType
    if ($code -ne $null) { $ws1.Cells.Item($row, 4).Value = $code }
    if ($codeSystem -ne $null) { $ws1.Cells.Item($row, 5).Value = $codeSystem }
    $formula = 'CONCATENATE("INSERT INTO dbo.attributes (id, name, value_type, code, code_system) VALUES (", A' + $row + ', ", ''", B' + $row + ', "'', ''", C' + $row + ', "'', ", IF(D' + $row + ' = "", "NULL", CONCATENATE("''", D' + $row + ', "''")), ", ", IF(E' + $row + ' = "", "NULL", CONCATENATE("''", E' + $row + ', "''")), ")")'
    $ws1.Cells.Item($row, 7).Formula = "=" + $formula
}

Set-AttributeRow 129 128 "SNP allele" "short_text" $null $null
Set-AttributeRow 130 129 "Star allele" "short_text" $null $null

# Set names in the order that matches original string-pool allocation:
# Star allele should get a lower shared-string index than SNP allele.
$ws1.Cells.Item(130, 2).Value = "Star allele"
$ws1.Cells.Item(129, 2).Value = "SNP allele"
